# Week 40 profiles update
#
# Flip the "done" flag (column C) from 0 -> 1 for the four profiles that
# were completed this week: rows 39, 104, 106 and 108.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C39").Value  = 1
$ws.Range("C104").Value = 1
$ws.Range("C106").Value = 1
$ws.Range("C108").Value = 1

# Leave the view parked on the last row touched, matching where the author's
# cursor ended up (C106) when the workbook was saved.
$ws.Range("C106").Select()
